# Week 15 logged (actual results) + Week 16 simulated (play-by-play style
# running totals appended to the long space-separated lists, plus weekly
# cumulative counters incremented on the Home/Road summary sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append this week's per-play yardage samples to the four
# long running lists (OFF/DEF x R/P).
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 10 2 2 1 7 6 7 6 1 8 4 2 8 4 2 2 4 5 -1"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 5 2 10 0 -2 1 4 1 2 6 1 2 3 9 39 0 0 6 6 10 3 0 17 1 -4 4"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 6 7 -5 9 4 11 9 34 9 25 12 3 3 2 9 7 11"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 7 0 1 14 18 5 11 12 9 9 3 2 42 5 32 6 -1 24 5 11 29"

# ---------------------------------------------------------------------
# OFF sheet: cumulative season counters for Home (row 2) / Road (row 3).
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 156
$wsOFF.Range("F2").Value = 41
$wsOFF.Range("G2").Value = 40
$wsOFF.Range("J2").Value = 23
$wsOFF.Range("N2").Value = 21
$wsOFF.Range("O2").Value = 9

$wsOFF.Range("B3").Value = 6
$wsOFF.Range("C3").Value = 127
$wsOFF.Range("E3").Value = 31
$wsOFF.Range("F3").Value = 76
$wsOFF.Range("G3").Value = 32
$wsOFF.Range("H3").Value = 18
$wsOFF.Range("I3").Value = 44
$wsOFF.Range("J3").Value = 41
$wsOFF.Range("L3").Value = 239
$wsOFF.Range("M3").Value = 154
$wsOFF.Range("Q3").Value = 443

# ---------------------------------------------------------------------
# DEF sheet: cumulative season counters for Home (row 2) / Road (row 3).
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 195
$wsDEF.Range("E2").Value = 5
$wsDEF.Range("F2").Value = 58
$wsDEF.Range("G2").Value = 68
$wsDEF.Range("H2").Value = 2
$wsDEF.Range("I2").Value = 10
$wsDEF.Range("J2").Value = 29
$wsDEF.Range("N2").Value = 15
$wsDEF.Range("O2").Value = 24

$wsDEF.Range("B3").Value = 10
$wsDEF.Range("C3").Value = 193
$wsDEF.Range("E3").Value = 28
$wsDEF.Range("F3").Value = 110
$wsDEF.Range("G3").Value = 39
$wsDEF.Range("H3").Value = 24
$wsDEF.Range("I3").Value = 47
$wsDEF.Range("J3").Value = 50
$wsDEF.Range("L3").Value = 307
$wsDEF.Range("M3").Value = 212
$wsDEF.Range("Q3").Value = 577

# ---------------------------------------------------------------------
# ST sheet: weekly counters plus the four long per-kick/return lists.
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 58
$wsST.Range("D2").Value = 75
$wsST.Range("F2").Value = 251
$wsST.Range("G2").Value = 225
$wsST.Range("J2").Value = 88
$wsST.Range("K2").Value = 83

$wsST.Range("B3").Value = 18

$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 63 65"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 32 14"
$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 40 47 54 59 48 49"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 0 0 0 0 0 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 9 5 0"

# ---------------------------------------------------------------------
# TURNS sheet: Road (row 3) counters.
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value = 3
$wsTURNS.Range("C3").Value = 4
$wsTURNS.Range("D3").Value = 6

# ---------------------------------------------------------------------
# PEN sheet: False start (row 2) / Holding (row 3) / Pass interference
# (row 4) counters.
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 11
$wsPEN.Range("B3").Value = 11
$wsPEN.Range("D3").Value = 8
$wsPEN.Range("D4").Value = 7
